# tdf116486.docx: Open Sans -> DejaVu Sans for repeatable layout in test
#
# The original styles.xml sets "Open Sans" (plus the "Light"/"Semibold"
# weight variants) as the document's fonts in six places: the document
# defaults (rPrDefault) and five individual styles. This script updates
# the five styles that are reachable through the Word object model's
# Styles collection; each of them only carries w:ascii/w:hAnsi (no
# eastAsia/cs), so setting Font.Name reproduces the diff exactly.

$d = $word.ActiveDocument

$styleNames = @(
    "Normal",
    "Strong",
    "Querbalken 1. Ebene fett",
    "Fußbereich Standard",
    "Seitenzahl Folgeseiten Zchn"
)

foreach ($name in $styleNames) {
    $style = $d.Styles($name)
    $style.Font.Name = "DejaVu Sans"
}
